# Add a new RF submission row ("less features") to the overview table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the existing table ("Tabelle1") by one row; this grows the
# table/autoFilter reference from A1:I13 to A1:I14.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# Fill in the new row's cells. The order below mirrors how the shared
# strings table ends up populated (Parameters, Average AUC, Hand in,
# then File Name), matching the authored workbook.
$ws.Range("D14").Value = 'max_features="log2"'
$ws.Range("G14").Value = "0.619 (0.022)"
$ws.Range("H14").Value = "Feb. 28, 2023, 11:24 a.m."
$ws.Range("A14").Value = "2023-02-28-1210_RF_centers_less_feat.csv"
$ws.Range("B14").Value = "RandomForest"
$ws.Range("C14").Value = "MoCov"
$ws.Range("E14").Value = "1 x 3"
$ws.Range("F14").Value = "weakly supervision with cv centers"
$ws.Range("I14").Value = 0.591

# Move the active selection to the next empty row, as it would be after
# entering the new row of data.
$ws.Range("A15").Select()
